$d = $word.ActiveDocument

$old = "COLPAERT; HUBERTUS. Metalografia dos produtos siderúrgicos comuns, 3ª Edição, Editora Edgard Blücher Ltda, SãoPaulo – 1974.COUTINHO, TELMO DE AZEVEDO. Metalografia de Não-Ferrosos, Editora Edgard Blücher Ltda, São Paulo – 1980.PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.REED-HILL, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982.Nondestructive Characterization of Materials. Series. Plenum Press, New York.YACOBI, B.G.; HOLT, D.B.; KAZMERSKI, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994."

$new = "COLPAERT; HUBERTUS. Metalografia dos produtos siderúrgicos comuns, 3ª Edição, Editora Edgard Blücher Ltda, São^lPaulo – 1974.^lCOUTINHO, TELMO DE AZEVEDO. Metalografia de Não-Ferrosos, Editora Edgard Blücher Ltda, São Paulo – 1980.^lPADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.^lMURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.^lWU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.^lREED-HILL, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982.^lNondestructive Characterization of Materials. Series. Plenum Press, New York.^lYACOBI, B.G.; HOLT, D.B.; KAZMERSKI, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994."

$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Replacement.ClearFormatting()
$range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
